$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert a new row at row 5; this shifts the old rows 5-15 down to 6-16
# and automatically updates dimension + the total formula (SUM(D2:D14) -> SUM(D2:D15)).
$ws.Rows.Item(5).Insert()

# --- Row 3: "Brainstorm 10 use cases" is now done by Naman and marked Completed ---
$ws.Range("B3").Value = "Naman"
$ws.Range("C3").Value = "Completed"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

# --- Row 4: "Collate the best use cases" is now done by Hunter and marked Completed ---
$ws.Range("A4").Value = "Collate the best use cases"
$ws.Range("B4").Value = "Hunter"
$ws.Range("C4").Value = "Completed"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# --- Row 5 (newly inserted): UC1 Start Game, Stephen, Not Started, 1 day ---
$ws.Range("A5").Value = "UC1: Start Game"
$ws.Range("B5").Value = "Stephen"
$ws.Range("C5").Value = "Not Started"
$ws.Range("D5").Value = 1

# --- Row 6: UC2 Create Character, Hunter, Not Started, 1 day ---
$ws.Range("A6").Value = "UC2: Create Character"
$ws.Range("B6").Value = "Hunter"
$ws.Range("D6").Value = 1

# --- Row 7: UC5 Buy New Goods, Naman, Not Started, 1 day ---
$ws.Range("A7").Value = "UC5: Buy New Goods"
$ws.Range("B7").Value = "Naman"
$ws.Range("D7").Value = 1

# --- Row 9: UC9 Leave Project, Bhavesh, Not Started, 1 day ---
$ws.Range("A9").Value = "UC9: Leave Project"
$ws.Range("B9").Value = "Bhavesh"
$ws.Range("D9").Value = 1

# --- Row 8: UC6 Upgrade/Repair Ship, Pranil, Not Started, 1 day ---
$ws.Range("A8").Value = "UC6: Upgrade/Repair Ship"
$ws.Range("B8").Value = "Pranil"
$ws.Range("D8").Value = 1

# Update the selected cell to match the final state of the workbook.
$ws.Range("D5").Select()
